$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-12 (columns A-E): text, x, y, width, height
$data = @(
    @("Greek ", 152, 644, 56, 23),
    @("mythology ", 208, 644, 95, 23),
    @("gift ", 859, 926.8, 33, 23),
    @("of ", 892, 926.8, 22, 23),
    @("prophecy. ", 152, 951.8, 89, 23),
    @("Trojan ", 348, 1001.8, 62, 23),
    @("Horse ", 410, 1001.8, 57, 23),
    @("trick, ", 467, 1001.8, 49, 23),
    @("Agamemnon’s ", 222, 1026.8, 126, 23),
    @("Bronze ", 406, 1193.2, 65, 23),
    @("Age. ", 471, 1193.2, 42, 23)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
